# DOMA-2542 Localization for Excel template (ticket_report_status_executor)
#
# The "{d.tickets[i + 1].<field>}" placeholders (third data row of the
# template) need their spacing normalised to "{d.tickets[i+1].<field>}" so
# the template engine's expression parser accepts them the same way it does
# "{d.tickets[i].<field>}".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rows = $used.Rows.Count
$cols = $used.Columns.Count

for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $text = $cell.Text
        if ($text.Contains("tickets[i + 1]")) {
            $cell.Value = $text.Replace("tickets[i + 1]", "tickets[i+1]")
        }
    }
}
